$d = $word.ActiveDocument

$d.Content.Find.Execute("on Windows Vista or Windows 7", $true, $false, $false, $false, $false,
                         $true, 1, $false, "on Windows", 2)
